$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# Update row 3 values
$ws.Range("D3").Value = "abcd"
$ws.Range("C3").Value = "ABCD"

# Clear E3 and F3 entirely (delete cell contents)
$ws.Range("E3:F3").ClearContents()

# Update selection to F3
$ws.Range("F3").Select()
